$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "Fe3" row) entirely; this shifts row 4 up to row 3
$ws.Rows.Item(2).Delete()

# Update the value in B3 (previously B4) from 170.556 to 3.021
$ws.Range("B3").Value = 3.021
